$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(86, 1).NumberFormat = "@"
$ws.Cells.Item(86, 1).Value = "2025-01-23"
$ws.Cells.Item(86, 1).Style = "Normal"

$ws.Cells.Item(86, 2).NumberFormat = "@"
$ws.Cells.Item(86, 2).Value = "42.6"
$ws.Cells.Item(86, 2).Style = "Normal"
